# Apply cell-value corrections to the stock report (qty/value recalculations
# and a few row-pair swaps) as described by the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 80: column(s) F,G
$ws.Cells.Item(80, 6).Value = 20
$ws.Cells.Item(80, 7).Value = 4601.6

# Row 81: column(s) F,G
$ws.Cells.Item(81, 6).Value = 19
$ws.Cells.Item(81, 7).Value = 10867.43

# Row 98: column(s) F,G
$ws.Cells.Item(98, 6).Value = 18
$ws.Cells.Item(98, 7).Value = 2318.4

# Row 102: column(s) B
$ws.Cells.Item(102, 2).Value = 191624.61

# Row 125: column(s) F,G
$ws.Cells.Item(125, 6).Value = 202
$ws.Cells.Item(125, 7).Value = 22660.36

# Row 147: column(s) B
$ws.Cells.Item(147, 2).Value = 108291.13

# Row 151: column(s) B,F,G
$ws.Cells.Item(151, 2).Value = 65258
$ws.Cells.Item(151, 6).Value = 2
$ws.Cells.Item(151, 7).Value = 64287.16

# Row 152: column(s) B,F,G
$ws.Cells.Item(152, 2).Value = 64196
$ws.Cells.Item(152, 6).Value = 1
$ws.Cells.Item(152, 7).Value = 32143.58

# Row 178: column(s) F,G
$ws.Cells.Item(178, 6).Value = 113
$ws.Cells.Item(178, 7).Value = 5882.78

# Row 184: column(s) B
$ws.Cells.Item(184, 2).Value = 34687.23

# Row 235: column(s) F,G
$ws.Cells.Item(235, 6).Value = 74
$ws.Cells.Item(235, 7).Value = 6625.96

# Row 238: column(s) F,G
$ws.Cells.Item(238, 6).Value = 112
$ws.Cells.Item(238, 7).Value = 9184

# Row 251: column(s) B
$ws.Cells.Item(251, 2).Value = 114230.47

# Row 267: column(s) F,G
$ws.Cells.Item(267, 6).Value = 132
$ws.Cells.Item(267, 7).Value = 8553.6

# Row 269: column(s) F,G
$ws.Cells.Item(269, 6).Value = 27
$ws.Cells.Item(269, 7).Value = 2345.76

# Row 274: column(s) B
$ws.Cells.Item(274, 2).Value = 12545.52

# Row 279: column(s) F,G
$ws.Cells.Item(279, 6).Value = 272
$ws.Cells.Item(279, 7).Value = 34456.96

# Row 283: column(s) F,G
$ws.Cells.Item(283, 6).Value = 136
$ws.Cells.Item(283, 7).Value = 10104.8

# Row 284: column(s) B
$ws.Cells.Item(284, 2).Value = 130450.55

# Row 298: column(s) F,G
$ws.Cells.Item(298, 6).Value = 40
$ws.Cells.Item(298, 7).Value = 4585.2

# Row 307: column(s) B
$ws.Cells.Item(307, 2).Value = 18048.65

# Row 315: column(s) F,G
$ws.Cells.Item(315, 6).Value = 5
$ws.Cells.Item(315, 7).Value = 27980.3

# Row 319: column(s) B
$ws.Cells.Item(319, 2).Value = 65301.03

# Row 340: column(s) F,G
$ws.Cells.Item(340, 6).Value = 231
$ws.Cells.Item(340, 7).Value = 9812.879999999999

# Row 342: column(s) F,G
$ws.Cells.Item(342, 6).Value = 51
$ws.Cells.Item(342, 7).Value = 6487.71

# Row 359: column(s) F,G
$ws.Cells.Item(359, 6).Value = 39
$ws.Cells.Item(359, 7).Value = 3786.9

# Row 363: column(s) F,G
$ws.Cells.Item(363, 6).Value = 250
$ws.Cells.Item(363, 7).Value = 11717.5

# Row 376: column(s) B
$ws.Cells.Item(376, 2).Value = 202631.56

# Row 391: column(s) B,E,F,G
$ws.Cells.Item(391, 2).Value = 55356
$ws.Cells.Item(391, 5).Value = 54.04
$ws.Cells.Item(391, 6).Value = -158
$ws.Cells.Item(391, 7).Value = -7527.12

# Row 392: column(s) B,E,F,G
$ws.Cells.Item(392, 2).Value = 63510
$ws.Cells.Item(392, 5).Value = 50.66
$ws.Cells.Item(392, 6).Value = 88
$ws.Cells.Item(392, 7).Value = 4192.32

# Row 401: column(s) B,E,F,G
$ws.Cells.Item(401, 2).Value = 60325
$ws.Cells.Item(401, 5).Value = 151.57
$ws.Cells.Item(401, 6).Value = -102
$ws.Cells.Item(401, 7).Value = -12939.72

# Row 402: column(s) B,E,F,G
$ws.Cells.Item(402, 2).Value = 63560
$ws.Cells.Item(402, 5).Value = 134.87
$ws.Cells.Item(402, 6).Value = 1
$ws.Cells.Item(402, 7).Value = 126.86

# Row 409: column(s) F,G
$ws.Cells.Item(409, 6).Value = 267
$ws.Cells.Item(409, 7).Value = 45745.11

# Row 413: column(s) B
$ws.Cells.Item(413, 2).Value = 71803.82000000001

# Row 460: column(s) F,G
$ws.Cells.Item(460, 6).Value = 471
$ws.Cells.Item(460, 7).Value = 66217.89

# Row 462: column(s) B
$ws.Cells.Item(462, 2).Value = 145885.74

# Row 472: column(s) F,G
$ws.Cells.Item(472, 6).Value = 392
$ws.Cells.Item(472, 7).Value = 65068.08

# Row 474: column(s) B
$ws.Cells.Item(474, 2).Value = 114250.09

# Row 479: column(s) F,G
$ws.Cells.Item(479, 6).Value = 68
$ws.Cells.Item(479, 7).Value = 65947.08

# Row 480: column(s) B
$ws.Cells.Item(480, 2).Value = 65947.08

# Row 490: column(s) F,G
$ws.Cells.Item(490, 6).Value = 719
$ws.Cells.Item(490, 7).Value = 69455.39999999999

# Row 493: column(s) B
$ws.Cells.Item(493, 2).Value = 88751.46000000001

# Row 563: column(s) B,E,F,G
$ws.Cells.Item(563, 2).Value = 45718
$ws.Cells.Item(563, 5).Value = 19.38
$ws.Cells.Item(563, 6).Value = -294
$ws.Cells.Item(563, 7).Value = -4768.68

# Row 564: column(s) B,E,F,G
$ws.Cells.Item(564, 2).Value = 64927
$ws.Cells.Item(564, 5).Value = 17.26
$ws.Cells.Item(564, 6).Value = 106
$ws.Cells.Item(564, 7).Value = 1719.32

# Row 568: column(s) B,E,F,G
$ws.Cells.Item(568, 2).Value = 45709
$ws.Cells.Item(568, 5).Value = 15.69
$ws.Cells.Item(568, 6).Value = -300
$ws.Cells.Item(568, 7).Value = -3945

# Row 569: column(s) B,E,F,G
$ws.Cells.Item(569, 2).Value = 64925
$ws.Cells.Item(569, 5).Value = 13.97
$ws.Cells.Item(569, 6).Value = 111
$ws.Cells.Item(569, 7).Value = 1459.65

# Row 570: column(s) B,E,F,G
$ws.Cells.Item(570, 2).Value = 64919
$ws.Cells.Item(570, 5).Value = 27.97
$ws.Cells.Item(570, 6).Value = 61
$ws.Cells.Item(570, 7).Value = 1604.3

# Row 571: column(s) B,E,F,G
$ws.Cells.Item(571, 2).Value = 45702
$ws.Cells.Item(571, 5).Value = 31.43
$ws.Cells.Item(571, 6).Value = -215
$ws.Cells.Item(571, 7).Value = -5654.5

# Row 573: column(s) B,E,F,G
$ws.Cells.Item(573, 2).Value = 65067
$ws.Cells.Item(573, 5).Value = 15.65
$ws.Cells.Item(573, 6).Value = 126
$ws.Cells.Item(573, 7).Value = 1855.98

# Row 574: column(s) B,E,F,G
$ws.Cells.Item(574, 2).Value = 53595
$ws.Cells.Item(574, 5).Value = 17.61
$ws.Cells.Item(574, 6).Value = -335
$ws.Cells.Item(574, 7).Value = -4934.55

# Row 613: column(s) F,G
$ws.Cells.Item(613, 6).Value = 95
$ws.Cells.Item(613, 7).Value = 21111.85

# Row 616: column(s) B
$ws.Cells.Item(616, 2).Value = 168228.61

# Row 636: column(s) F,G
$ws.Cells.Item(636, 6).Value = 23
$ws.Cells.Item(636, 7).Value = 3572.36

# Row 637: column(s) B
$ws.Cells.Item(637, 2).Value = 80559.67999999999

# Row 673: column(s) B,E,F,G
$ws.Cells.Item(673, 2).Value = 64830
$ws.Cells.Item(673, 5).Value = 34.9
$ws.Cells.Item(673, 6).Value = 92
$ws.Cells.Item(673, 7).Value = 3020.36

# Row 674: column(s) B,E,F,G
$ws.Cells.Item(674, 2).Value = 60022
$ws.Cells.Item(674, 5).Value = 37.22
$ws.Cells.Item(674, 6).Value = -113
$ws.Cells.Item(674, 7).Value = -3709.79

# Row 680: column(s) F,G
$ws.Cells.Item(680, 6).Value = 294
$ws.Cells.Item(680, 7).Value = 47181.12

# Row 696: column(s) B
$ws.Cells.Item(696, 2).Value = 220730.14

# Row 707: column(s) F,G
$ws.Cells.Item(707, 6).Value = 92
$ws.Cells.Item(707, 7).Value = 2767.36

# Row 709: column(s) B
$ws.Cells.Item(709, 2).Value = 47297.09

# Row 740: column(s) F,G
$ws.Cells.Item(740, 6).Value = 348
$ws.Cells.Item(740, 7).Value = 42403.8

# Row 744: column(s) F,G
$ws.Cells.Item(744, 6).Value = 13
$ws.Cells.Item(744, 7).Value = 1077.57

# Row 746: column(s) F,G
$ws.Cells.Item(746, 6).Value = 20
$ws.Cells.Item(746, 7).Value = 1657.8

# Row 747: column(s) B
$ws.Cells.Item(747, 2).Value = 57345.15

# Row 802: column(s) F,G
$ws.Cells.Item(802, 6).Value = 13
$ws.Cells.Item(802, 7).Value = 1060.28

# Row 804: column(s) F,G
$ws.Cells.Item(804, 6).Value = 295
$ws.Cells.Item(804, 7).Value = 39264.5

# Row 806: column(s) F,G
$ws.Cells.Item(806, 6).Value = 94
$ws.Cells.Item(806, 7).Value = 11346.74

# Row 807: column(s) B
$ws.Cells.Item(807, 2).Value = 54245.28

# Row 822: column(s) F,G
$ws.Cells.Item(822, 6).Value = 109
$ws.Cells.Item(822, 7).Value = 15561.93

# Row 830: column(s) F,G
$ws.Cells.Item(830, 6).Value = 473
$ws.Cells.Item(830, 7).Value = 37173.07

# Row 843: column(s) B
$ws.Cells.Item(843, 2).Value = 402382.01

# Row 884: column(s) F,G
$ws.Cells.Item(884, 6).Value = 122
$ws.Cells.Item(884, 7).Value = 9797.82

# Row 891: column(s) B
$ws.Cells.Item(891, 2).Value = 30757.37

# Row 896: column(s) F,G
$ws.Cells.Item(896, 6).Value = 447
$ws.Cells.Item(896, 7).Value = 13512.81

# Row 897: column(s) F,G
$ws.Cells.Item(897, 6).Value = 2060
$ws.Cells.Item(897, 7).Value = 336006.6

# Row 899: column(s) F,G
$ws.Cells.Item(899, 6).Value = 133
$ws.Cells.Item(899, 7).Value = 37621.71

# Row 903: column(s) B
$ws.Cells.Item(903, 2).Value = 402942.08

# Row 915: column(s) F,G
$ws.Cells.Item(915, 6).Value = 55
$ws.Cells.Item(915, 7).Value = 8753.799999999999

# Row 919: column(s) B
$ws.Cells.Item(919, 2).Value = 24593.42

# Row 948: column(s) B
$ws.Cells.Item(948, 2).Value = 6150790.89

# Row 949: column(s) B
$ws.Cells.Item(949, 2).Value = 6150790.89
